# "added latest grid editor version"
#
# - begroting!B19 moves from the "Neutral" highlight to the "Good" highlight
#   (and the now-unused "Neutral" cell style is dropped from the workbook)
# - two more weeks of hours-spent data are appended (week 16 gets an hours
#   value, week 17 is added) which bumps the totals/budget formulas
# - the view state changes: Sheet1 (the "Sheet1" tab) becomes the active /
#   selected sheet, with a new selection on each sheet

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "begroting"
$ws2 = $wb.Worksheets.Item(2)   # "Sheet1"

# --- begroting!B19: Neutral -> Good -------------------------------------
$ws1.Range("B19").Style = "Good"

# the "Neutral" cell style is no longer used anywhere, drop it
$wb.Styles.Item("Neutral").Delete()

# --- new hours-spent rows -------------------------------------------------
# existing row 39 (week 16) gains an hours value, and a brand new row is
# inserted below it for week 17 - this pushes the old rows 40/41 ("budget"
# label + "Total hours spent so far" totals) down by one row
$ws1.Rows.Item(40).Insert()

$ws1.Range("B39").Value = 30

$ws1.Range("A40").Value = 17
$ws1.Range("B40").Value = 30

# the totals formulas keep referencing the same (now one-row-lower) range
$ws1.Range("B43").Formula = "=SUM(B27:B40)"
$ws1.Range("C43").Formula = "=B43*135"

# --- selection / active sheet ---------------------------------------------
$ws1.Range("B40").Select()
$ws2.Range("F20").Select()

$ws2.Activate()
